# Apply the commit's edits to the workbook.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# CX (sheet1): condense the data down to two "jobs" (rows 2-3), keep the
# two marker rows (4-5), add a "SecondaryPP" column G, and append an
# "endofsheet" marker far below the data (row 30).
# ---------------------------------------------------------------------
$cx = $wb.Worksheets.Item("CX")

# Drop the old row 6 (blank D="03000" spacer) and the old rows 2-4
# (first job block) - this leaves the former rows 5,7,8,9 as the new
# rows 2,3,4,5.
$cx.Rows("6:6").Delete() | Out-Null
$cx.Rows("2:4").Delete() | Out-Null

# Row 2 (former row 5) keeps its ClampUnit ("0350") but gets a new CAD
# name and PrimaryPlast value.
$cx.Range("A2").Value = "cx_0250_emptyy.asm"
$cx.Range("D2").Value = "01000"

# New column header.
$cx.Range("G1").Value = "SecondaryPP"

# Extend the row styling into the new column G for the two data rows and
# the trailing blank row.
$cx.Range("F2").Copy() | Out-Null
$cx.Range("G2").PasteSpecial(-4122) | Out-Null
$cx.Range("F3").Copy() | Out-Null
$cx.Range("G3").PasteSpecial(-4122) | Out-Null
$cx.Range("F5").Copy() | Out-Null
$cx.Range("G5").PasteSpecial(-4122) | Out-Null

# Marker row far below the data (style copied from the D4 marker cell).
$cx.Range("D4").Copy() | Out-Null
$cx.Range("A30").PasteSpecial(-4122) | Out-Null
$cx.Range("A30").Value = "endofsheet"

$cx.Activate()
$cx.Range("A30").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1

# ---------------------------------------------------------------------
# GX (sheet2): becomes the (slightly tweaked) former CX content, and
# gets an "END_OF_WORKSHEET" marker appended at row 25. It ends up being
# the active tab.
# ---------------------------------------------------------------------
$gx = $wb.Worksheets.Item("GX")

$gx.Range("A2").Value = "cx_0250_gm1811.asm"

$gx.Range("D2").Copy() | Out-Null
$gx.Range("A25").PasteSpecial(-4122) | Out-Null
$gx.Range("A25").Value = "END_OF_WORKSHEET"

$gx.Activate()
$gx.Range("A25").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1

# ---------------------------------------------------------------------
# PX (sheet3): only the saved selection changes.
# ---------------------------------------------------------------------
$px = $wb.Worksheets.Item("PX")
$px.Range("B14").Select() | Out-Null

# ---------------------------------------------------------------------
# GXH (sheet7): no longer the active tab; its view just scrolls down.
# ---------------------------------------------------------------------
$gxh = $wb.Worksheets.Item("GXH")
$gxh.Range("A4").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1

# GX is the sheet left active/visible when the workbook was saved.
$gx.Activate()
